# Add the new data row (row 15) that the app appended to the "lista" sheet.
# The values are stored as plain text (matching the app's inline-string
# export), so we prefix them with a single-quote to force text entry
# instead of letting Excel auto-convert the numeric-looking strings, then
# clear the resulting cell format so no stray number-format/style is left
# behind on the new cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "'7847.0"
$ws.Range("B15").Value = "'7417.0"
$ws.Range("C15").Value = "'17"

$ws.Range("A15:C15").ClearFormats()
